$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values for rows 2-5
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 51

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 25

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 20

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3

# Remove row 6 entirely (data no longer present, dimension shrinks to A1:B5)
$ws.Range("A6:B6").Delete()
